$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# Row 17
$ws.Range("H17").Value = 2526.6304
$ws.Range("J17").Value = 2625.568
$ws.Range("L17").Value = 7876.704000000001
$ws.Range("N17").Value = -8212.704000000002

# Row 112
$ws.Range("H112").Value = 1447.1578
$ws.Range("J112").Value = 1546.8235
$ws.Range("L112").Value = 4640.470499999999
$ws.Range("N112").Value = -6856.470499999999

# Row 137
$ws.Range("H137").Value = 1440.5758
$ws.Range("I137").Value = 893.5769
$ws.Range("J137").Value = 3472.2856
$ws.Range("K137").Value = 2680.7307
$ws.Range("L137").Value = 10416.8568
$ws.Range("M137").Value = -130.7307000000001
$ws.Range("N137").Value = -15516.8568

# Row 138
$ws.Range("H138").Value = 1731.678
$ws.Range("I138").Value = 599.05
$ws.Range("J138").Value = 4116.1577
$ws.Range("K138").Value = 1797.15
$ws.Range("L138").Value = 12348.4731
$ws.Range("M138").Value = 3342.85
$ws.Range("N138").Value = -22628.4731

$ws = $wb.Worksheets.Item("ARM")

# Row 61
$ws.Range("H61").Value = 1357.8334
$ws.Range("I61").Value = 1310.3636
$ws.Range("J61").Value = 1398
$ws.Range("K61").Value = 1310.3636
$ws.Range("L61").Value = 1398
$ws.Range("M61").Value = -1098.3636
$ws.Range("N61").Value = -1822

# Row 74
$ws.Range("H74").Value = 4387280.5
$ws.Range("I74").Value = 6579712
$ws.Range("J74").Value = 2417.4736
$ws.Range("K74").Value = 6579712
$ws.Range("L74").Value = 2417.4736
$ws.Range("M74").Value = -6578838
$ws.Range("N74").Value = -4165.473599999999

# Row 77
$ws.Range("H77").Value = 4387280.5
$ws.Range("I77").Value = 6579712
$ws.Range("J77").Value = 2417.4736
$ws.Range("K77").Value = 32898560
$ws.Range("L77").Value = 12087.368
$ws.Range("M77").Value = -32894192
$ws.Range("N77").Value = -20823.368

# Row 136
$ws.Range("H136").Value = 1357.8334
$ws.Range("I136").Value = 1310.3636
$ws.Range("J136").Value = 1398
$ws.Range("K136").Value = 3931.0908
$ws.Range("L136").Value = 4194
$ws.Range("M136").Value = -1381.0908
$ws.Range("N136").Value = -9294

$ws = $wb.Worksheets.Item("BSM")

# Row 134
$ws.Range("H134").Value = 1326.95
$ws.Range("I134").Value = 1338.4615
$ws.Range("J134").Value = 1305.5714
$ws.Range("K134").Value = 4015.3845
$ws.Range("L134").Value = 3916.7142
$ws.Range("M134").Value = -1480.3845
$ws.Range("N134").Value = -8986.7142

$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 4904164.5
$ws.Range("I31").Value = 7445565
$ws.Range("J31").Value = 2892.0715
$ws.Range("K31").Value = 7445565
$ws.Range("L31").Value = 2892.0715
$ws.Range("M31").Value = -7445270
$ws.Range("N31").Value = -3482.0715

# Row 34
$ws.Range("H34").Value = 4904164.5
$ws.Range("I34").Value = 7445565
$ws.Range("J34").Value = 2892.0715
$ws.Range("K34").Value = 7445565
$ws.Range("L34").Value = 2892.0715
$ws.Range("M34").Value = -7445363
$ws.Range("N34").Value = -3296.0715

# Row 58
$ws.Range("H58").Value = 1009.40424
$ws.Range("I58").Value = 473.73077
$ws.Range("J58").Value = 1672.619
$ws.Range("K58").Value = 473.73077
$ws.Range("L58").Value = 1672.619
$ws.Range("M58").Value = -270.73077
$ws.Range("N58").Value = -2078.619

# Row 132
$ws.Range("H132").Value = 1546.6
$ws.Range("I132").Value = 1469.5385
$ws.Range("J132").Value = 1689.7142
$ws.Range("K132").Value = 4408.6155
$ws.Range("L132").Value = 5069.142599999999
$ws.Range("M132").Value = -1878.6155
$ws.Range("N132").Value = -10129.1426

# Row 134
$ws.Range("H134").Value = 5222.533
$ws.Range("I134").Value = 8732
$ws.Range("J134").Value = 2151.75
$ws.Range("K134").Value = 26196
$ws.Range("L134").Value = 6455.25
$ws.Range("M134").Value = -23661
$ws.Range("N134").Value = -11525.25

# Row 136
$ws.Range("H136").Value = 1009.40424
$ws.Range("I136").Value = 473.73077
$ws.Range("J136").Value = 1672.619
$ws.Range("K136").Value = 1421.19231
$ws.Range("L136").Value = 5017.857
$ws.Range("M136").Value = 1128.80769
$ws.Range("N136").Value = -10117.857

$ws = $wb.Worksheets.Item("CUL")

# Row 113
$ws.Range("H113").Value = 2275
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2275
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 6825
$ws.Range("N113").Value = -11165
$ws.Range("M113").ClearContents()

# Row 129
$ws.Range("H129").Value = 1078.4615
$ws.Range("I129").Value = 1078.4615
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 3235.3845
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = 1764.6155
$ws.Range("N129").ClearContents()

# Row 130
$ws.Range("H130").Value = 984.3
$ws.Range("I130").Value = 720.4286
$ws.Range("K130").Value = 2161.2858
$ws.Range("M130").Value = 2858.7142

# Row 131
$ws.Range("H131").Value = 944.3778
$ws.Range("I131").Value = 370.92307
$ws.Range("J131").Value = 1177.3438
$ws.Range("K131").Value = 1112.76921
$ws.Range("L131").Value = 3532.0314
$ws.Range("M131").Value = 3927.23079
$ws.Range("N131").Value = -13612.0314

# Row 136
$ws.Range("H136").Value = 700.44183
$ws.Range("I136").Value = 483.96667
$ws.Range("K136").Value = 1451.90001
$ws.Range("M136").Value = 3648.09999

# Row 139
$ws.Range("H139").Value = 3092.7917
$ws.Range("I139").Value = 1444.8823
$ws.Range("J139").Value = 7094.857
$ws.Range("K139").Value = 4334.6469
$ws.Range("L139").Value = 21284.571
$ws.Range("M139").Value = 805.3531000000003
$ws.Range("N139").Value = -31564.571

# Row 140
$ws.Range("H140").Value = 2330.9092
$ws.Range("I140").Value = 1664
$ws.Range("J140").Value = 9000
$ws.Range("K140").Value = 4992
$ws.Range("L140").Value = 27000
$ws.Range("M140").Value = 188
$ws.Range("N140").Value = -37360

$ws = $wb.Worksheets.Item("LTW")

# Row 68
$ws.Range("H68").Value = 13535362
$ws.Range("I68").Value = 26026788
$ws.Range("J68").Value = 2985.1667
$ws.Range("K68").Value = 26026788
$ws.Range("L68").Value = 2985.1667
$ws.Range("M68").Value = -26026039
$ws.Range("N68").Value = -4483.1667

# Row 71
$ws.Range("H71").Value = 13535362
$ws.Range("I71").Value = 26026788
$ws.Range("J71").Value = 2985.1667
$ws.Range("K71").Value = 130133940
$ws.Range("L71").Value = 14925.8335
$ws.Range("M71").Value = -130130196
$ws.Range("N71").Value = -22413.8335

# Row 132
$ws.Range("H132").Value = 7026019.5
$ws.Range("I132").Value = 13894159
$ws.Range("J132").Value = 1786.7273
$ws.Range("K132").Value = 41682477
$ws.Range("L132").Value = 5360.1819
$ws.Range("M132").Value = -41679947
$ws.Range("N132").Value = -10420.1819

# Row 136
$ws.Range("H136").Value = 2653.8193
$ws.Range("I136").Value = 3306.7874
$ws.Range("J136").Value = 1426.24
$ws.Range("K136").Value = 9920.3622
$ws.Range("L136").Value = 4278.72
$ws.Range("M136").Value = -7370.3622
$ws.Range("N136").Value = -9378.720000000001

$ws = $wb.Worksheets.Item("WVR")

# Row 132
$ws.Range("H132").Value = 1425.8889
$ws.Range("I132").Value = 753.93475
$ws.Range("J132").Value = 3244.1177
$ws.Range("K132").Value = 2261.80425
$ws.Range("L132").Value = 9732.3531
$ws.Range("M132").Value = 268.1957499999999
$ws.Range("N132").Value = -14792.3531

# Row 136
$ws.Range("H136").Value = 480.66
$ws.Range("I136").Value = 268.2
$ws.Range("J136").Value = 976.4
$ws.Range("K136").Value = 804.5999999999999
$ws.Range("L136").Value = 2929.2
$ws.Range("M136").Value = 1745.4
$ws.Range("N136").Value = -8029.2
